$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2710
$ws.Range("I76").Value = 2622.8125
$ws.Range("K76").Value = 2622.8125
$ws.Range("M76").Value = -2307.8125

$ws.Range("H79").Value = 2710
$ws.Range("I79").Value = 2622.8125
$ws.Range("K79").Value = 2622.8125
$ws.Range("M79").Value = -1530.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H74").Value = 1987.762
$ws.Range("I74").Value = 1395.625
$ws.Range("J74").Value = 3882.6
$ws.Range("K74").Value = 1395.625
$ws.Range("L74").Value = 3882.6
$ws.Range("M74").Value = -521.625
$ws.Range("N74").Value = -5630.6

$ws.Range("H77").Value = 1987.762
$ws.Range("I77").Value = 1395.625
$ws.Range("J77").Value = 3882.6
$ws.Range("K77").Value = 6978.125
$ws.Range("L77").Value = 19413
$ws.Range("M77").Value = -2610.125
$ws.Range("N77").Value = -28149

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H39").Value = 9842.166999999999
$ws.Range("J39").Value = 9842.166999999999
$ws.Range("L39").Value = 9842.166999999999
$ws.Range("N39").Value = -10620.167

$ws.Range("H56").Value = 16110
$ws.Range("J56").Value = 16110
$ws.Range("L56").Value = 16110
$ws.Range("N56").Value = -17588

$ws.Range("H105").Value = 1835.6154
$ws.Range("I105").Value = 1803.5
$ws.Range("J105").Value = 1887
$ws.Range("K105").Value = 1803.5
$ws.Range("L105").Value = 1887
$ws.Range("M105").Value = -56.5
$ws.Range("N105").Value = -5381

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1548.2727
$ws.Range("I16").Value = 628.875
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 628.875
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -341.875
$ws.Range("N16").Value = -4574

$ws.Range("H58").Value = 8623791
$ws.Range("I58").Value = 1545.3939
$ws.Range("J58").Value = 20005156
$ws.Range("K58").Value = 1545.3939
$ws.Range("L58").Value = 20005156
$ws.Range("M58").Value = -1342.3939
$ws.Range("N58").Value = -20005562

$ws.Range("H113").Value = 1548.2727
$ws.Range("I113").Value = 628.875
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 628.875
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 1541.125
$ws.Range("N113").Value = -8340

$ws.Range("H136").Value = 8623791
$ws.Range("I136").Value = 1545.3939
$ws.Range("J136").Value = 20005156
$ws.Range("K136").Value = 4636.1817
$ws.Range("L136").Value = 60015468
$ws.Range("M136").Value = -2086.1817
$ws.Range("N136").Value = -60020568

$ws.Range("H140").Value = 17125.9
$ws.Range("J140").Value = 17125.9
$ws.Range("L140").Value = 17125.9
$ws.Range("N140").Value = -27485.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 174.45454
$ws.Range("J33").Value = 217.33333
$ws.Range("L33").Value = 1303.99998
$ws.Range("N33").Value = -1869.99998

$ws.Range("H122").Value = 1500.4286
$ws.Range("I122").Value = 700.8
$ws.Range("K122").Value = 6307.2
$ws.Range("M122").Value = -3857.2

$ws.Range("H131").Value = 1427.5264
$ws.Range("I131").Value = 4659
$ws.Range("J131").Value = 1116.8077
$ws.Range("K131").Value = 13977
$ws.Range("L131").Value = 3350.4231
$ws.Range("M131").Value = -8937
$ws.Range("N131").Value = -13430.4231

$ws.Range("H132").Value = 2630.15
$ws.Range("I132").Value = 1600.7273
$ws.Range("J132").Value = 3888.3333
$ws.Range("K132").Value = 14406.5457
$ws.Range("L132").Value = 34994.9997
$ws.Range("M132").Value = -11876.5457
$ws.Range("N132").Value = -40054.9997

$ws.Range("H134").Value = 2937.1177
$ws.Range("I134").Value = 1622.7273
$ws.Range("J134").Value = 5346.8335
$ws.Range("K134").Value = 4868.1819
$ws.Range("L134").Value = 16040.5005
$ws.Range("M134").Value = 201.8181000000004
$ws.Range("N134").Value = -26180.5005

$ws.Range("H136").Value = 1821.7826
$ws.Range("I136").Value = 1409.3334
$ws.Range("J136").Value = 3306.6
$ws.Range("K136").Value = 4228.0002
$ws.Range("L136").Value = 9919.799999999999
$ws.Range("M136").Value = 871.9997999999996
$ws.Range("N136").Value = -20119.8

$ws.Range("H137").Value = 3685.6086
$ws.Range("I137").Value = 2781.25
$ws.Range("J137").Value = 5752.7144
$ws.Range("K137").Value = 8343.75
$ws.Range("L137").Value = 17258.1432
$ws.Range("M137").Value = -3243.75
$ws.Range("N137").Value = -27458.1432

$ws.Range("H138").Value = 1300.6923
$ws.Range("I138").Value = 1075.75
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 3227.25
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = 1912.75
$ws.Range("N138").Value = -22280

$ws.Range("H139").Value = 10003106
$ws.Range("I139").Value = 12502882
$ws.Range("J139").Value = 3998
$ws.Range("K139").Value = 37508646
$ws.Range("L139").Value = 11994
$ws.Range("M139").Value = -37503506
$ws.Range("N139").Value = -22274

$ws.Range("H140").Value = 5382419.5
$ws.Range("I140").Value = 11906091
$ws.Range("J140").Value = 9984.235000000001
$ws.Range("K140").Value = 35718273
$ws.Range("L140").Value = 29952.705
$ws.Range("M140").Value = -35713093
$ws.Range("N140").Value = -40312.705

$ws.Range("H141").Value = 2995
$ws.Range("J141").Value = 4990
$ws.Range("L141").Value = 14970
$ws.Range("N141").Value = -25330

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4875
$ws.Range("I70").Value = 5071.4287
$ws.Range("J70").Value = 4600
$ws.Range("K70").Value = 5071.4287
$ws.Range("L70").Value = 4600
$ws.Range("M70").Value = -4801.4287
$ws.Range("N70").Value = -5140

$ws.Range("H73").Value = 4875
$ws.Range("I73").Value = 5071.4287
$ws.Range("J73").Value = 4600
$ws.Range("K73").Value = 5071.4287
$ws.Range("L73").Value = 4600
$ws.Range("M73").Value = -4135.4287
$ws.Range("N73").Value = -6472

$ws.Range("H80").Value = 2553.875
$ws.Range("I80").Value = 2737.4375
$ws.Range("J80").Value = 2186.75
$ws.Range("K80").Value = 2737.4375
$ws.Range("L80").Value = 2186.75
$ws.Range("M80").Value = -1739.4375
$ws.Range("N80").Value = -4182.75

$ws.Range("H83").Value = 2553.875
$ws.Range("I83").Value = 2737.4375
$ws.Range("J83").Value = 2186.75
$ws.Range("K83").Value = 13687.1875
$ws.Range("L83").Value = 10933.75
$ws.Range("M83").Value = -8695.1875
$ws.Range("N83").Value = -20917.75

$ws.Range("H103").Value = 24734.166
$ws.Range("I103").Value = 14000
$ws.Range("J103").Value = 30101.25
$ws.Range("K103").Value = 14000
$ws.Range("L103").Value = 30101.25
$ws.Range("M103").Value = -12828
$ws.Range("N103").Value = -32445.25

$ws.Range("H108").Value = 28000
$ws.Range("J108").Value = 28000
$ws.Range("L108").Value = 28000
$ws.Range("N108").Value = -35680

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2480
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 3866.6667
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 3866.6667
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -4456.6667

$ws.Range("H27").Value = 2480
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 3866.6667
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 3866.6667
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -4080.6667

$ws.Range("H61").Value = 90912584
$ws.Range("I61").Value = 125000700
$ws.Range("K61").Value = 125000700
$ws.Range("M61").Value = -125000498

$ws.Range("H113").Value = 90912584
$ws.Range("I113").Value = 125000700
$ws.Range("K113").Value = 125000700
$ws.Range("M113").Value = -124998530

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3500
$ws.Range("I113").Value = 200
$ws.Range("J113").Value = 5150
$ws.Range("K113").Value = 600
$ws.Range("L113").Value = 15450
$ws.Range("M113").Value = 1570
$ws.Range("N113").Value = -19790

$ws.Range("H136").Value = 2525.6843
$ws.Range("I136").Value = 2020.04
$ws.Range("J136").Value = 3498.077
$ws.Range("K136").Value = 6060.12
$ws.Range("L136").Value = 10494.231
$ws.Range("M136").Value = -3510.12
$ws.Range("N136").Value = -15594.231
